$d = $word.ActiveDocument

function Get-ParaIndexAtPos($doc, $pos) {
    return $doc.Range(0, $pos + 1).Paragraphs.Count
}

# ---------------------------------------------------------------------
# 1. Delete the whole "License Information" (Heading2) paragraph.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'License Information' paragraph"
}
$d.Range($rng.Start, $rng.End + 1).Delete()

# ---------------------------------------------------------------------
# 2. Locate the paragraph that starts with the bold run
#    "Pertanyaan Terjemahan (unfoldingWord)" (the license-credit
#    paragraph) and merge it with the very next paragraph
#    ("This PDF version is provided under the same license.") by
#    deleting the paragraph mark between them.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("is based on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the license-credit paragraph"
}
$creditIdx = Get-ParaIndexAtPos $d $rng2.Start
$creditPara = $d.Paragraphs.Item($creditIdx)

$mark = $d.Range($creditPara.Range.End - 1, $creditPara.Range.End)
$mark.Delete()

# Re-fetch the (now merged) paragraph and wipe its text, leaving the
# paragraph mark (and therefore the leading/trailing empty runs) intact.
$merged = $d.Paragraphs.Item($creditIdx)
$clearRange = $d.Range($merged.Range.Start, $merged.Range.End - 1)
$clearRange.Text = ""

# ---------------------------------------------------------------------
# 3. Rebuild the paragraph content run by run, in order, using the
#    collapsed insertion-point pattern so each chunk becomes its own run
#    with independent character formatting.
# ---------------------------------------------------------------------
$merged2 = $d.Paragraphs.Item($creditIdx)
$pos = $merged2.Range.Start

$ip = $d.Range($pos, $pos)
$ip.InsertBefore("unfoldingWord® Translation Questions")
$ip.Font.Bold = 1

$ip = $d.Range($ip.End, $ip.End)
$ip.InsertBefore(" © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. ")
$ip.Font.Bold = 0

$ip = $d.Range($ip.End, $ip.End)
$ip.InsertBefore("unfoldingWord® Translation Questions")
$ip.Font.Bold = 0

$ip = $d.Range($ip.End, $ip.End)
$ip.InsertBefore(" has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from ")
$ip.Font.Bold = 0

$ip = $d.Range($ip.End, $ip.End)
$ip.InsertBefore("unfoldingWord® Translation Questions")
$ip.Font.Bold = 0

$ip = $d.Range($ip.End, $ip.End)
$ip.InsertBefore(" © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual")
$ip.Font.Bold = 0
